$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("F5").Value = "（環境業務課）087-834-0389、087-861-4524（環境指導課）087-839-2380（適正処理対策室）087-839-2370"
$ws.Range("G5").Value = "（環境業務課）http://www.city.takamatsu.kagawa.jp/kurashi/shinotorikumi/soshikihyo/kankyogyomu.html（環境指導課）http://www.city.takamatsu.kagawa.jp/kurashi/shinotorikumi/soshikihyo/kankyoshido.html（適正処理対策室）http://www.city.takamatsu.kagawa.jp/kurashi/shinotorikumi/soshikihyo/tekiseisyoritaisaku.html"
$ws.Range("K5").Value = "16:45"
$ws.Range("L5").Value = "利用可能曜日は祝祭日及び年末年始を除く。"
